# chore: update Sheets via scheduled runner
# Refreshes market-board-derived profit figures (columns H:N, i.e.
# currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on a handful of
# leve rows across several class sheets. A couple of rows also gain/lose a
# cell entirely where the refreshed source no longer (or now does) produce
# an HQ profit figure.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 177.5
$ws.Range("I4").Value = 177.5
$ws.Range("K4").Value = 177.5
$ws.Range("M4").Value = -63.5
# Row 33
$ws.Range("H33").Value = 5996.8887
$ws.Range("I33").Value = 7874.154
$ws.Range("J33").Value = 1116
$ws.Range("K33").Value = 7874.154
$ws.Range("L33").Value = 1116
$ws.Range("M33").Value = -7645.154
$ws.Range("N33").Value = -1574
# Row 99
$ws.Range("H99").Value = 10000280
$ws.Range("I99").Value = 10000280
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 30000840
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -29999342
$ws.Range("N99").Value = $null
# Row 107
$ws.Range("H107").Value = 562.3182
$ws.Range("I107").Value = 613.05554
$ws.Range("J107").Value = 334
$ws.Range("K107").Value = 613.05554
$ws.Range("L107").Value = 334
$ws.Range("M107").Value = 1306.94446
$ws.Range("N107").Value = -4174
# Row 129
$ws.Range("H129").Value = 1020.5
$ws.Range("I129").Value = 454.55554
$ws.Range("J129").Value = 1151.1025
$ws.Range("K129").Value = 1363.66662
$ws.Range("L129").Value = 3453.3075
$ws.Range("M129").Value = 3636.33338
$ws.Range("N129").Value = -13453.3075

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 8315.666999999999
$ws.Range("I28").Value = 3260.2856
$ws.Range("K28").Value = 3260.2856
$ws.Range("M28").Value = -3068.2856
# Row 74
$ws.Range("H74").Value = 1238.2449
$ws.Range("I74").Value = 1279.7368
$ws.Range("J74").Value = 1094.909
$ws.Range("K74").Value = 1279.7368
$ws.Range("L74").Value = 1094.909
$ws.Range("M74").Value = -405.7367999999999
$ws.Range("N74").Value = -2842.909
# Row 77
$ws.Range("H77").Value = 1238.2449
$ws.Range("I77").Value = 1279.7368
$ws.Range("J77").Value = 1094.909
$ws.Range("K77").Value = 6398.683999999999
$ws.Range("L77").Value = 5474.545
$ws.Range("M77").Value = -2030.683999999999
$ws.Range("N77").Value = -14210.545
# Row 99
$ws.Range("H99").Value = 8315.666999999999
$ws.Range("I99").Value = 3260.2856
$ws.Range("K99").Value = 3260.2856
$ws.Range("M99").Value = -265.2856000000002
# Row 120
$ws.Range("H120").Value = 37990
$ws.Range("J120").Value = 37990
$ws.Range("L120").Value = 37990
$ws.Range("N120").Value = -47666
# Row 133
$ws.Range("H133").Value = 36899.8
$ws.Range("J133").Value = 36899.8
$ws.Range("L133").Value = 36899.8
$ws.Range("N133").Value = -41959.8
# Row 135
$ws.Range("H135").Value = 68857.25
$ws.Range("I135").Value = 35000
$ws.Range("J135").Value = 80143
$ws.Range("K135").Value = 35000
$ws.Range("L135").Value = 80143
$ws.Range("M135").Value = -29930
$ws.Range("N135").Value = -90283

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 49042.855
$ws.Range("J50").Value = 49042.855
$ws.Range("L50").Value = 49042.855
$ws.Range("N50").Value = -50292.855
# Row 51
$ws.Range("H51").Value = 44412.5
$ws.Range("J51").Value = 44412.5
$ws.Range("L51").Value = 44412.5
$ws.Range("N51").Value = -45884.5
# Row 58
$ws.Range("H58").Value = 1293.4
$ws.Range("I58").Value = 1221.6
$ws.Range("J58").Value = 1508.8
$ws.Range("K58").Value = 1221.6
$ws.Range("L58").Value = 1508.8
$ws.Range("M58").Value = -1018.6
$ws.Range("N58").Value = -1914.8
# Row 60
$ws.Range("H60").Value = 23500
$ws.Range("J60").Value = 30500
$ws.Range("L60").Value = 30500
$ws.Range("N60").Value = -31522
# Row 61
$ws.Range("H61").Value = 44412.5
$ws.Range("J61").Value = 44412.5
$ws.Range("L61").Value = 44412.5
$ws.Range("N61").Value = -45108.5
# Row 136
$ws.Range("H136").Value = 1293.4
$ws.Range("I136").Value = 1221.6
$ws.Range("J136").Value = 1508.8
$ws.Range("K136").Value = 3664.8
$ws.Range("L136").Value = 4526.4
$ws.Range("M136").Value = -1114.8
$ws.Range("N136").Value = -9626.4

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 274920.12
$ws.Range("I5").Value = 302.66666
$ws.Range("J5").Value = 475316.7
$ws.Range("K5").Value = 907.9999799999999
$ws.Range("L5").Value = 1425950.1
$ws.Range("M5").Value = -795.9999799999999
$ws.Range("N5").Value = -1426174.1
# Row 102
$ws.Range("H102").Value = 6163
$ws.Range("J102").Value = 6423.1875
$ws.Range("L102").Value = 19269.5625
$ws.Range("N102").Value = -24137.5625
# Row 122
$ws.Range("H122").Value = 520.5294
$ws.Range("I122").Value = 338.7857
$ws.Range("J122").Value = 1368.6666
$ws.Range("K122").Value = 3049.0713
$ws.Range("L122").Value = 12317.9994
$ws.Range("M122").Value = -599.0713000000001
$ws.Range("N122").Value = -17217.9994
# Row 132
$ws.Range("H132").Value = 1473.0769
$ws.Range("I132").Value = 1590
$ws.Range("K132").Value = 14310
$ws.Range("M132").Value = -11780
# Row 135
$ws.Range("H135").Value = 274920.12
$ws.Range("I135").Value = 302.66666
$ws.Range("J135").Value = 475316.7
$ws.Range("K135").Value = 2723.99994
$ws.Range("L135").Value = 4277850.3
$ws.Range("M135").Value = -188.9999399999997
$ws.Range("N135").Value = -4282920.3

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null
# Row 141
$ws.Range("H141").Value = 53554.25
$ws.Range("J141").Value = 53554.25
$ws.Range("L141").Value = 53554.25
$ws.Range("N141").Value = -63914.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 22068.889
$ws.Range("I40").Value = 27764.45
$ws.Range("J40").Value = 5795.857
$ws.Range("K40").Value = 27764.45
$ws.Range("L40").Value = 5795.857
$ws.Range("M40").Value = -27628.45
$ws.Range("N40").Value = -6067.857
# Row 136
$ws.Range("H136").Value = 7248634
$ws.Range("I136").Value = 2089.8572
$ws.Range("K136").Value = 6269.571599999999
$ws.Range("M136").Value = -3719.571599999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 34992
$ws.Range("I75").Value = 15000
$ws.Range("K75").Value = 15000
$ws.Range("M75").Value = -14064
# Row 78
$ws.Range("H78").Value = 34992
$ws.Range("I78").Value = 15000
$ws.Range("K78").Value = 45000
$ws.Range("M78").Value = -40320
# Row 136
$ws.Range("H136").Value = 5131523.5
$ws.Range("I136").Value = 11905452
$ws.Range("J136").Value = 5307.6484
$ws.Range("K136").Value = 35716356
$ws.Range("L136").Value = 15922.9452
$ws.Range("M136").Value = -35713806
$ws.Range("N136").Value = -21022.9452
